# Add methodology section in Chapter 3
# Fill in the two new "Actions taken" comments for the Chapter-3 methodology
# correction (rows 5 and 6, column C), shown in the thesis correction list
# in red text to flag them as newly-added responses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 / column C: pointer to where it was addressed in the thesis.
# (Set first so it claims the earlier shared-string slot, matching the
# original author's edit order.)
$c6 = $ws.Range("C6")
$c6.Value = "To explain this more in Section 3.1"
$c6.Font.Color = 255          # pure red (RGB 255,0,0)
$c6.HorizontalAlignment = -4108 # xlCenter
$c6.VerticalAlignment = -4108   # xlCenter
$c6.WrapText = $true

# Row 5 / column C: what was added (kept vertically centered, no wrap).
$c5 = $ws.Range("C5")
$c5.Value = "Add a section. Make a figure of procedures. Include a photo of the setup. "
$c5.Font.Color = 255          # pure red (RGB 255,0,0)
$c5.VerticalAlignment = -4108 # xlCenter

# Reflect the author's cursor ending up on C8 after reviewing the new text.
$ws.Range("C8").Select()
